$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 52, shifting the existing weekly records
# (rows 52-77) down to rows 53-78, and add the new week's Mango price
# record for "Feria Lagunitas de Puerto Montt" in the freshly opened row.
$ws.Rows(52).Insert()

$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 44466
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100108
$ws.Cells.Item(52, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(52, 9).Value = 100108002
$ws.Cells.Item(52, 10).Value = "Mango"
$ws.Cells.Item(52, 11).Value = "Sin especificar"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 80
$ws.Cells.Item(52, 14).Value = 10000
$ws.Cells.Item(52, 15).Value = 10000
$ws.Cells.Item(52, 16).Value = 10000
$ws.Cells.Item(52, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(52, 18).Value = "Brasil"
$ws.Cells.Item(52, 19).Value = 2500
$ws.Cells.Item(52, 20).Value = 4
